$wb = $excel.ActiveWorkbook

# --- Sheet "Настройки": append new settings rows 16-34 ---
$wsSettings = $wb.Worksheets.Item("Настройки")
$wsSettings.Range("A16").Value = "IS_Service_type"
$wsSettings.Range("B16").Value = "Тип сервиса (ИСУ, КИС, ЛИС, ПУ, ..)"
$wsSettings.Range("C16").Value = "AX"
$wsSettings.Range("D16").Value = "AY"

$wsSettings.Range("A17").Value = "IS_Product_type"
$wsSettings.Range("B17").Value = "Тип системы (SAP, БК, ЛИМС, MES,…)"
$wsSettings.Range("C17").Value = "AZ"
$wsSettings.Range("D17").Value = "BA"

$wsSettings.Range("A18").Value = "Pdr_Proj"
$wsSettings.Range("B18").Value = "Группировка Подразделение+Проект"
$wsSettings.Range("C18").Value = "BB"
$wsSettings.Range("D18").Value = "BC"

$wsSettings.Range("A19").Value = "Proj_Pdr"
$wsSettings.Range("B19").Value = "Группировка Проект+Подразделение"
$wsSettings.Range("C19").Value = "BD"
$wsSettings.Range("D19").Value = "BE"

$wsSettings.Range("A20").Value = "Portfolio"
$wsSettings.Range("B20").Value = "Портфель проектов"
$wsSettings.Range("C20").Value = "AV"
$wsSettings.Range("D20").Value = "AW"

$wsSettings.Range("A21").Value = "Personal_email"
$wsSettings.Range("B21").Value = "Признак отправлять сообщение лично или в общей массе"
$wsSettings.Range("C21").Value = "AR"

$wsSettings.Range("A22").Value = "user_email"
$wsSettings.Range("B22").Value = "Почтовый адрес пользователя"
$wsSettings.Range("C22").Value = "AS"

$wsSettings.Range("A23").Value = "boss_email"
$wsSettings.Range("B23").Value = "Почтовый адрес руководителя данного пользователя"
$wsSettings.Range("C23").Value = "AT"

$wsSettings.Range("A24").Value = "Contract"
$wsSettings.Range("B24").Value = "Доходный договор"
$wsSettings.Range("C24").Value = "AU"

$wsSettings.Range("A25").Value = "FN"
$wsSettings.Range("B25").Value = "Функциональное направление (или подразделение)"
$wsSettings.Range("C25").Value = "C"
$wsSettings.Range("D25").Value = "BF"

$wsSettings.Range("A26").Value = "UHCost_KV1"
$wsSettings.Range("B26").Value = "Часовая ставка в 1-м квартале"
$wsSettings.Range("C26").Value = "BF"

$wsSettings.Range("A27").Value = "UMCost_KV1"
$wsSettings.Range("B27").Value = "Месячная ставка в 1-м квартале"
$wsSettings.Range("C27").Value = "BG"

$wsSettings.Range("A28").Value = "UHCost_KV2"
$wsSettings.Range("B28").Value = "Часовая ставка во 2-м квартале"
$wsSettings.Range("C28").Value = "BH"

$wsSettings.Range("A29").Value = "UMCost_KV2"
$wsSettings.Range("B29").Value = "Месячная ставка во 2-м квартале"
$wsSettings.Range("C29").Value = "BI"

$wsSettings.Range("A30").Value = "UHCost_KV3"
$wsSettings.Range("B30").Value = "Часовая ставка в 3-м квартале"
$wsSettings.Range("C30").Value = "BJ"

$wsSettings.Range("A31").Value = "UMCost_KV3"
$wsSettings.Range("B31").Value = "Месячная ставка в 3-м квартале"
$wsSettings.Range("C31").Value = "BK"

$wsSettings.Range("A32").Value = "UHCost_KV4"
$wsSettings.Range("B32").Value = "Часовая ставка в 4-м квартале"
$wsSettings.Range("C32").Value = "BL"

$wsSettings.Range("A33").Value = "UMCost_KV4"
$wsSettings.Range("B33").Value = "Месячная ставка в 4-м квартале"
$wsSettings.Range("C33").Value = "BM"

$wsSettings.Range("A34").Value = "ISDogName"
$wsSettings.Range("B34").Value = "Название ИС из договора"
$wsSettings.Range("C34").Value = "BO"

# --- Sheet "ИсходныеДанные": append new header columns AX1:BO1 ---
$wsSource = $wb.Worksheets.Item("ИсходныеДанные")
$wsSource.Range("AX1").Value = "IS_Service_type"
$wsSource.Range("AY1").Value = "IS_Service_type_Month"
$wsSource.Range("AZ1").Value = "IS_Product_type"
$wsSource.Range("BA1").Value = "IS_Product_type_Month"
$wsSource.Range("BB1").Value = "Pdr_Proj"
$wsSource.Range("BC1").Value = "Pdr_Proj_Month"
$wsSource.Range("BD1").Value = "Proj_Pdr"
$wsSource.Range("BE1").Value = "Proj_Pdr_Month"
$wsSource.Range("BF1").Value = "FN_Month"
$wsSource.Range("BG1").Value = "UHCost_KV1"
$wsSource.Range("BH1").Value = "UMCost_KV1"
$wsSource.Range("BI1").Value = "UHCost_KV2"
$wsSource.Range("BJ1").Value = "UMCost_KV2"
$wsSource.Range("BK1").Value = "UHCost_KV3"
$wsSource.Range("BL1").Value = "UMCost_KV3"
$wsSource.Range("BM1").Value = "UHCost_KV4"
$wsSource.Range("BN1").Value = "UMCost_KV4"
$wsSource.Range("BO1").Value = "ISDogName"

# --- Restore selection on "Настройки" sheet to match target, then refocus original active sheet ---
$wsSettings.Range("A1:D34").Select()
$wsReport = $wb.Worksheets.Item("Отчет")
$wsReport.Activate()
